$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "NAICS Code:" header label in B1 was renamed to "Codes:"
$ws.Range("B1").Value = "Codes:"

# Reflect the author's new active cell selection (B1) on the sheet view
$ws.Range("B1").Select()
